# Visitors sheet update: new "created At3" text timestamps replace the old
# numeric "updated" / "created At" / "updated At" date columns (I, J, K),
# and two freshly-submitted visitor rows (6 and 7) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the old "created At" / "updated At" columns (J:K). ------------
# Column I ("updated") becomes the sole, renamed "created At3" column and
# inherits the MM/DD/YYYY HH:MM:SS format that used to live on J/K.
$ws.Range("J1:K5").EntireColumn.Delete()

# --- 2. Re-purpose column I as a text "created At3" column. ----------------
$ws.Range("I1:I7").NumberFormat = "MM/DD/YYYY HH:MM:SS"

$ws.Range("I1").Value = "created At3"

# Existing rows 2 & 3 never got a "created At3" stamp - clear the old date
# value but keep the cell (and its format) in place.
$ws.Range("I2").ClearContents()
$ws.Range("I3").ClearContents()

# Existing rows 4 & 5 get the human-readable timestamp text instead of the
# old serial-date value.
$ws.Range("I4").Value = "Sat Jul 30 2022 18:31:07 GMT+0300 (Eastern European Summer Time)3"
$ws.Range("I5").Value = "Sat Jul 30 2022 18:32:19 GMT+0300 (Eastern European Summer Time)3"

# --- 3. Append the two new visitor records (rows 6 & 7). -------------------
$ws.Range("A6").Value = '"62e66c45b301c57149bdf071"'
$ws.Range("B6").Value = "mouaz test"
$ws.Range("C6").Value = "123123123"
$ws.Range("D6").Value = $false
$ws.Range("G6").Value = "google Search"
$ws.Range("H6").Value = "now@hotmail.com"
$ws.Range("I6").Value = "Sun Jul 31 2022 14:49:25 GMT+0300 (Eastern European Summer Time)3"

$ws.Range("A7").Value = '"62e6709542df06499605257d"'
$ws.Range("B7").Value = "mouaz last test time"
$ws.Range("C7").Value = "123123123"
$ws.Range("D7").Value = $false
$ws.Range("G7").Value = "instagram"
$ws.Range("H7").Value = "now@hotmail.com"
$ws.Range("I7").Value = "Sun Jul 31 2022 15:07:49 GMT+0300 (Eastern European Summer Time)3"
